$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 10070
$ws.Range("J17").Value = 10070
$ws.Range("L17").Value = 30210
$ws.Range("N17").Value = -30546

$ws.Range("H76").Value = 4239
$ws.Range("J76").Value = 2795
$ws.Range("L76").Value = 2795
$ws.Range("N76").Value = -3425

$ws.Range("H79").Value = 4239
$ws.Range("J79").Value = 2795
$ws.Range("L79").Value = 2795
$ws.Range("N79").Value = -4979

$ws.Range("H80").Value = 3276.5
$ws.Range("I80").Value = 6325.7144
$ws.Range("J80").Value = 904.8889
$ws.Range("K80").Value = 18977.1432
$ws.Range("L80").Value = 2714.6667
$ws.Range("M80").Value = -17979.1432
$ws.Range("N80").Value = -4710.6667

$ws.Range("H83").Value = 3276.5
$ws.Range("I83").Value = 6325.7144
$ws.Range("J83").Value = 904.8889
$ws.Range("K83").Value = 56931.4296
$ws.Range("L83").Value = 8144.0001
$ws.Range("M83").Value = -51939.4296
$ws.Range("N83").Value = -18128.0001

$ws.Range("H132").Value = 7253475
$ws.Range("I132").Value = 11116645
$ws.Range("J132").Value = 10031.5625
$ws.Range("K132").Value = 33349935
$ws.Range("L132").Value = 30094.6875
$ws.Range("M132").Value = -33347405
$ws.Range("N132").Value = -35154.6875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5317.544
$ws.Range("I32").Value = 4520.1113
$ws.Range("J32").Value = 19671.334
$ws.Range("K32").Value = 4520.1113
$ws.Range("L32").Value = 19671.334
$ws.Range("M32").Value = -4233.1113
$ws.Range("N32").Value = -20245.334

$ws.Range("N61").ClearContents()
$ws.Range("H61").Value = 71431064
$ws.Range("I61").Value = 71431064
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 71431064
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -71430852

$ws.Range("N136").ClearContents()
$ws.Range("H136").Value = 71431064
$ws.Range("I136").Value = 71431064
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 214293192
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -214290642

$ws.Range("H141").Value = 31891.857
$ws.Range("J141").Value = 31891.857
$ws.Range("L141").Value = 31891.857
$ws.Range("N141").Value = -42251.857

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 25000834
$ws.Range("I94").Value = 27778370
$ws.Range("J94").Value = 3010
$ws.Range("K94").Value = 27778370
$ws.Range("L94").Value = 3010
$ws.Range("M94").Value = -27777919
$ws.Range("N94").Value = -3912

$ws.Range("H107").Value = 1367.0714
$ws.Range("I107").Value = 1133.6
$ws.Range("K107").Value = 1133.6
$ws.Range("M107").Value = 786.4000000000001

$ws.Range("H130").Value = 32499.5
$ws.Range("J130").Value = 32499.5
$ws.Range("L130").Value = 32499.5
$ws.Range("N130").Value = -42539.5

$ws.Range("H134").Value = 7132.5264
$ws.Range("I134").Value = 1210.25
$ws.Range("J134").Value = 17285
$ws.Range("K134").Value = 3630.75
$ws.Range("L134").Value = 51855
$ws.Range("M134").Value = -1095.75
$ws.Range("N134").Value = -56925

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 55556772
$ws.Range("J16").Value = 1122.8572
$ws.Range("L16").Value = 1122.8572
$ws.Range("N16").Value = -1696.8572

$ws.Range("N43").ClearContents()
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0

$ws.Range("N101").ClearContents()
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0

$ws.Range("H113").Value = 55556772
$ws.Range("J113").Value = 1122.8572
$ws.Range("L113").Value = 1122.8572
$ws.Range("N113").Value = -5462.8572

$ws.Range("H140").Value = 37500
$ws.Range("J140").Value = 37500
$ws.Range("L140").Value = 37500
$ws.Range("N140").Value = -47860

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 1909.75
$ws.Range("J87").Value = 3700
$ws.Range("L87").Value = 11100
$ws.Range("N87").Value = -13596

$ws.Range("H90").Value = 1909.75
$ws.Range("J90").Value = 3700
$ws.Range("L90").Value = 33300
$ws.Range("N90").Value = -45780

$ws.Range("H131").Value = 18871180
$ws.Range("J131").Value = 3698.2827
$ws.Range("L131").Value = 11094.8481
$ws.Range("N131").Value = -21174.8481

$ws.Range("H139").Value = 2070.9656
$ws.Range("I139").Value = 2318.842
$ws.Range("K139").Value = 6956.526
$ws.Range("M139").Value = -1816.526

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3787.9333
$ws.Range("I70").Value = 3972.1667
$ws.Range("J70").Value = 3665.111
$ws.Range("K70").Value = 3972.1667
$ws.Range("L70").Value = 3665.111
$ws.Range("M70").Value = -3702.1667
$ws.Range("N70").Value = -4205.111

$ws.Range("H73").Value = 3787.9333
$ws.Range("I73").Value = 3972.1667
$ws.Range("J73").Value = 3665.111
$ws.Range("K73").Value = 3972.1667
$ws.Range("L73").Value = 3665.111
$ws.Range("M73").Value = -3036.1667
$ws.Range("N73").Value = -5537.111

$ws.Range("H86").Value = 25329.334
$ws.Range("J86").Value = 25329.334
$ws.Range("L86").Value = 25329.334
$ws.Range("N86").Value = -27701.334

$ws.Range("H89").Value = 25329.334
$ws.Range("J89").Value = 25329.334
$ws.Range("L89").Value = 75988.00199999999
$ws.Range("N89").Value = -87844.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 985.5714
$ws.Range("I22").Value = 724.25
$ws.Range("J22").Value = 1334
$ws.Range("K22").Value = 724.25
$ws.Range("L22").Value = 1334
$ws.Range("M22").Value = -429.25
$ws.Range("N22").Value = -1924

$ws.Range("H27").Value = 985.5714
$ws.Range("I27").Value = 724.25
$ws.Range("J27").Value = 1334
$ws.Range("K27").Value = 724.25
$ws.Range("L27").Value = 1334
$ws.Range("M27").Value = -617.25
$ws.Range("N27").Value = -1548

$ws.Range("H122").Value = 16677130
$ws.Range("I122").Value = 22737296
$ws.Range("J122").Value = 11676
$ws.Range("K122").Value = 68211888
$ws.Range("L122").Value = 35028
$ws.Range("M122").Value = -68209438
$ws.Range("N122").Value = -39928

$ws.Range("H132").Value = 84241.03999999999
$ws.Range("I132").Value = 20736.727
$ws.Range("J132").Value = 127900.25
$ws.Range("K132").Value = 62210.181
$ws.Range("L132").Value = 383700.75
$ws.Range("M132").Value = -59680.181
$ws.Range("N132").Value = -388760.75

$ws.Range("H136").Value = 11240.65
$ws.Range("I136").Value = 14154.333
$ws.Range("J136").Value = 2499.6
$ws.Range("K136").Value = 42462.999
$ws.Range("L136").Value = 7498.799999999999
$ws.Range("M136").Value = -39912.999
$ws.Range("N136").Value = -12598.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2327.476
$ws.Range("I132").Value = 1739.0834
$ws.Range("J132").Value = 3112
$ws.Range("K132").Value = 5217.2502
$ws.Range("L132").Value = 9336
$ws.Range("M132").Value = -2687.2502
$ws.Range("N132").Value = -14396

$ws.Range("H136").Value = 971.03125
$ws.Range("I136").Value = 861.4583
$ws.Range("K136").Value = 2584.3749
$ws.Range("M136").Value = -34.3748999999998
